$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component List")

# --- Row 54: the "Wideband connector" BOM line got re-sourced from the old
#     6-pos housing/receptacle to a Molex Mini-Fit Jr. 6-pos header, matching
#     the look of the other "...POS Header" rows (e.g. row 21).
#
# Copy formatting (number format / font / fill / borders) one cell at a time
# from row 21 onto the row 54 cells whose look changes. Values are set
# separately below, so only formatting is carried over.
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C54").PasteSpecial(-4122) | Out-Null

$ws.Range("F21").Copy() | Out-Null
$ws.Range("F54").PasteSpecial(-4122) | Out-Null

$ws.Range("G21").Copy() | Out-Null
$ws.Range("G54").PasteSpecial(-4122) | Out-Null

$ws.Range("H21").Copy() | Out-Null
$ws.Range("H54").PasteSpecial(-4122) | Out-Null

$ws.Range("I21").Copy() | Out-Null
$ws.Range("I54").PasteSpecial(-4122) | Out-Null

$ws.Range("L21").Copy() | Out-Null
$ws.Range("L54").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# New part values.
$ws.Range("B54").Value2 = "Wideband connector"
$ws.Range("C54").Value2 = "6 POS Header"
$ws.Range("D54").Value2 = "HEADER 6P MINIFIT"
$ws.Range("H54").Value2 = "39-30-1060"
$ws.Range("I54").Value2 = "WM1353-ND"
$ws.Range("J54").Value2 = "538-39-30-1060"
$ws.Range("K54").Value2 = 1
$ws.Range("L54").Value2 = 0.882

# Re-assert the row's formulas (values follow from the new inputs above).
$ws.Range("M54").Formula = "=K54*A54"
$ws.Range("N54").Formula = "=L54*A54"
$ws.Range("P54").Formula = '=IF(NOT(I54=""),A54&","&I54,"")'
$ws.Range("Q54").Formula = '=A54&"x "&C54'
$ws.Range("R54").Formula = '=IF(NOT(J54=""),J54&"|"&A54,"")'
$ws.Range("S54").Formula = '=H54&" "&A54'

# New column S got a dedicated width in the finished sheet (~16.125 chars;
# 15.33 is the nearest input that this engine's pixel-quantised column-width
# model resolves back to that same displayed width).
$ws.Columns("S").ColumnWidth = 15.33

# Leave the sheet scrolled/selected where the author finished editing.
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("S54").Select() | Out-Null
